$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ntng1"
$ws.Cells.Item(2, 3).Value = "Lrrc4c"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.015642
$ws.Cells.Item(2, 8).Value = 0.046926
$ws.Cells.Item(2, 9).Value = 0.08667416564617744
$ws.Cells.Item(2, 10).Value = 0.08667416564617747
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.1854593333333333
$ws.Cells.Item(2, 14).Value = 0.556378
$ws.Cells.Item(2, 15).Value = 0.2174741358149192
$ws.Cells.Item(2, 16).Value = 0.2174741358149192
$ws.Cells.Item(2, 17).Value = 0.002900954892
$ws.Cells.Item(2, 18).Value = 0.026108594028
$ws.Cells.Item(2, 19).Value = 0.0188493892713816
$ws.Cells.Item(2, 20).Value = 0.01884938927138161

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ntng1"
$ws.Cells.Item(3, 3).Value = "Lrrc4c"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.015642
$ws.Cells.Item(3, 8).Value = 0.046926
$ws.Cells.Item(3, 9).Value = 0.08667416564617744
$ws.Cells.Item(3, 10).Value = 0.08667416564617747
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.1535116666666667
$ws.Cells.Item(3, 14).Value = 0.460535
$ws.Cells.Item(3, 15).Value = 0.1800115229889101
$ws.Cells.Item(3, 16).Value = 0.1800115229889101
$ws.Cells.Item(3, 17).Value = 0.00240122949
$ws.Cells.Item(3, 18).Value = 0.02161106541
$ws.Cells.Item(3, 19).Value = 0.01560234856176147
$ws.Cells.Item(3, 20).Value = 0.01560234856176148

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Ntng1"
$ws.Cells.Item(4, 3).Value = "Lrrc4c"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.015642
$ws.Cells.Item(4, 8).Value = 0.046926
$ws.Cells.Item(4, 9).Value = 0.08667416564617744
$ws.Cells.Item(4, 10).Value = 0.08667416564617747
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.0009266666666666667
$ws.Cells.Item(4, 14).Value = 0.00278
$ws.Cells.Item(4, 15).Value = 0.001086631925715027
$ws.Cells.Item(4, 16).Value = 0.001086631925715027
$ws.Cells.Item(4, 17).Value = 0.00001449492
$ws.Cells.Item(4, 18).Value = 0.00013045428
$ws.Cells.Item(4, 19).Value = 0.00009418291552584905
$ws.Cells.Item(4, 20).Value = 0.00009418291552584908

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Ntng1"
$ws.Cells.Item(5, 3).Value = "Lrrc4c"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.015642
$ws.Cells.Item(5, 8).Value = 0.046926
$ws.Cells.Item(5, 9).Value = 0.08667416564617744
$ws.Cells.Item(5, 10).Value = 0.08667416564617747
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.5128903333333333
$ws.Cells.Item(5, 14).Value = 1.538671
$ws.Cells.Item(5, 15).Value = 0.6014277092704556
$ws.Cells.Item(5, 16).Value = 0.6014277092704556
$ws.Cells.Item(5, 17).Value = 0.008022630594
$ws.Cells.Item(5, 18).Value = 0.072203675346
$ws.Cells.Item(5, 19).Value = 0.05212824489750852
$ws.Cells.Item(5, 20).Value = 0.05212824489750854

# Row 6
$ws.Cells.Item(6, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(6, 2).Value = "Ntng1"
$ws.Cells.Item(6, 3).Value = "Lrrc4c"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.002436333333333333
$ws.Cells.Item(6, 8).Value = 0.007309
$ws.Cells.Item(6, 9).Value = 0.01350001015871609
$ws.Cells.Item(6, 10).Value = 0.01350001015871609
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.1854593333333333
$ws.Cells.Item(6, 14).Value = 0.556378
$ws.Cells.Item(6, 15).Value = 0.2174741358149192
$ws.Cells.Item(6, 16).Value = 0.2174741358149192
$ws.Cells.Item(6, 17).Value = 0.0004518407557777778
$ws.Cells.Item(6, 18).Value = 0.004066566802000001
$ws.Cells.Item(6, 19).Value = 0.002935903042759411
$ws.Cells.Item(6, 20).Value = 0.002935903042759412

# Row 7
$ws.Cells.Item(7, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(7, 2).Value = "Ntng1"
$ws.Cells.Item(7, 3).Value = "Lrrc4c"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.002436333333333333
$ws.Cells.Item(7, 8).Value = 0.007309
$ws.Cells.Item(7, 9).Value = 0.01350001015871609
$ws.Cells.Item(7, 10).Value = 0.01350001015871609
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.1535116666666667
$ws.Cells.Item(7, 14).Value = 0.460535
$ws.Cells.Item(7, 15).Value = 0.1800115229889101
$ws.Cells.Item(7, 16).Value = 0.1800115229889101
$ws.Cells.Item(7, 17).Value = 0.0003740055905555556
$ws.Cells.Item(7, 18).Value = 0.003366050315
$ws.Cells.Item(7, 19).Value = 0.00243015738903624
$ws.Cells.Item(7, 20).Value = 0.002430157389036241

# Row 8
$ws.Cells.Item(8, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 2).Value = "Ntng1"
$ws.Cells.Item(8, 3).Value = "Lrrc4c"
$ws.Cells.Item(8, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.002436333333333333
$ws.Cells.Item(8, 8).Value = 0.007309
$ws.Cells.Item(8, 9).Value = 0.01350001015871609
$ws.Cells.Item(8, 10).Value = 0.01350001015871609
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.0009266666666666667
$ws.Cells.Item(8, 14).Value = 0.00278
$ws.Cells.Item(8, 15).Value = 0.001086631925715027
$ws.Cells.Item(8, 16).Value = 0.001086631925715027
$ws.Cells.Item(8, 17).Value = 0.000002257668888888889
$ws.Cells.Item(8, 18).Value = 0.00002031902
$ws.Cells.Item(8, 19).Value = 0.00001466954203593809
$ws.Cells.Item(8, 20).Value = 0.00001466954203593809

# Row 9
$ws.Cells.Item(9, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 2).Value = "Ntng1"
$ws.Cells.Item(9, 3).Value = "Lrrc4c"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.002436333333333333
$ws.Cells.Item(9, 8).Value = 0.007309
$ws.Cells.Item(9, 9).Value = 0.01350001015871609
$ws.Cells.Item(9, 10).Value = 0.01350001015871609
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.5128903333333333
$ws.Cells.Item(9, 14).Value = 1.538671
$ws.Cells.Item(9, 15).Value = 0.6014277092704556
$ws.Cells.Item(9, 16).Value = 0.6014277092704556
$ws.Cells.Item(9, 17).Value = 0.001249571815444444
$ws.Cells.Item(9, 18).Value = 0.011246146339
$ws.Cells.Item(9, 19).Value = 0.008119280184884496
$ws.Cells.Item(9, 20).Value = 0.008119280184884496

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Ntng1"
$ws.Cells.Item(10, 3).Value = "Lrrc4c"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.1602236666666667
$ws.Cells.Item(10, 8).Value = 0.480671
$ws.Cells.Item(10, 9).Value = 0.8878182217813955
$ws.Cells.Item(10, 10).Value = 0.8878182217813956
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.1854593333333333
$ws.Cells.Item(10, 14).Value = 0.556378
$ws.Cells.Item(10, 15).Value = 0.2174741358149192
$ws.Cells.Item(10, 16).Value = 0.2174741358149192
$ws.Cells.Item(10, 17).Value = 0.02971497440422223
$ws.Cells.Item(10, 18).Value = 0.267434769638
$ws.Cells.Item(10, 19).Value = 0.1930775005426472
$ws.Cells.Item(10, 20).Value = 0.1930775005426473

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Ntng1"
$ws.Cells.Item(11, 3).Value = "Lrrc4c"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.1602236666666667
$ws.Cells.Item(11, 8).Value = 0.480671
$ws.Cells.Item(11, 9).Value = 0.8878182217813955
$ws.Cells.Item(11, 10).Value = 0.8878182217813956
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.1535116666666667
$ws.Cells.Item(11, 14).Value = 0.460535
$ws.Cells.Item(11, 15).Value = 0.1800115229889101
$ws.Cells.Item(11, 16).Value = 0.1800115229889101
$ws.Cells.Item(11, 17).Value = 0.02459620210944445
$ws.Cells.Item(11, 18).Value = 0.221365818985
$ws.Cells.Item(11, 19).Value = 0.159817510240175
$ws.Cells.Item(11, 20).Value = 0.159817510240175

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Ntng1"
$ws.Cells.Item(12, 3).Value = "Lrrc4c"
$ws.Cells.Item(12, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.1602236666666667
$ws.Cells.Item(12, 8).Value = 0.480671
$ws.Cells.Item(12, 9).Value = 0.8878182217813955
$ws.Cells.Item(12, 10).Value = 0.8878182217813956
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.0009266666666666667
$ws.Cells.Item(12, 14).Value = 0.00278
$ws.Cells.Item(12, 15).Value = 0.001086631925715027
$ws.Cells.Item(12, 16).Value = 0.001086631925715027
$ws.Cells.Item(12, 17).Value = 0.0001484739311111111
$ws.Cells.Item(12, 18).Value = 0.00133626538
$ws.Cells.Item(12, 19).Value = 0.0009647316240192088
$ws.Cells.Item(12, 20).Value = 0.0009647316240192089

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Ntng1"
$ws.Cells.Item(13, 3).Value = "Lrrc4c"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.1602236666666667
$ws.Cells.Item(13, 8).Value = 0.480671
$ws.Cells.Item(13, 9).Value = 0.8878182217813955
$ws.Cells.Item(13, 10).Value = 0.8878182217813956
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.5128903333333333
$ws.Cells.Item(13, 14).Value = 1.538671
$ws.Cells.Item(13, 15).Value = 0.6014277092704556
$ws.Cells.Item(13, 16).Value = 0.6014277092704556
$ws.Cells.Item(13, 17).Value = 0.08217716980455557
$ws.Cells.Item(13, 18).Value = 0.739594528241
$ws.Cells.Item(13, 19).Value = 0.533958479374554
$ws.Cells.Item(13, 20).Value = 0.5339584793745541

# Row 14
$ws.Cells.Item(14, 1).Value = "Neutrophils"
$ws.Cells.Item(14, 2).Value = "Ntng1"
$ws.Cells.Item(14, 3).Value = "Lrrc4c"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.002167
$ws.Cells.Item(14, 8).Value = 0.006501
$ws.Cells.Item(14, 9).Value = 0.01200760241371094
$ws.Cells.Item(14, 10).Value = 0.01200760241371094
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.1854593333333333
$ws.Cells.Item(14, 14).Value = 0.556378
$ws.Cells.Item(14, 15).Value = 0.2174741358149192
$ws.Cells.Item(14, 16).Value = 0.2174741358149192
$ws.Cells.Item(14, 17).Value = 0.0004018903753333333
$ws.Cells.Item(14, 18).Value = 0.003617013378
$ws.Cells.Item(14, 19).Value = 0.002611342958130925
$ws.Cells.Item(14, 20).Value = 0.002611342958130926

# Row 15
$ws.Cells.Item(15, 1).Value = "Neutrophils"
$ws.Cells.Item(15, 2).Value = "Ntng1"
$ws.Cells.Item(15, 3).Value = "Lrrc4c"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.002167
$ws.Cells.Item(15, 8).Value = 0.006501
$ws.Cells.Item(15, 9).Value = 0.01200760241371094
$ws.Cells.Item(15, 10).Value = 0.01200760241371094
$ws.Cells.Item(15, 11).Value = 2
$ws.Cells.Item(15, 12).Value = 0.6666666666666666
$ws.Cells.Item(15, 13).Value = 0.1535116666666667
$ws.Cells.Item(15, 14).Value = 0.460535
$ws.Cells.Item(15, 15).Value = 0.1800115229889101
$ws.Cells.Item(15, 16).Value = 0.1800115229889101
$ws.Cells.Item(15, 17).Value = 0.0003326597816666667
$ws.Cells.Item(15, 18).Value = 0.002993938035
$ws.Cells.Item(15, 19).Value = 0.002161506797937419
$ws.Cells.Item(15, 20).Value = 0.00216150679793742

# Row 16
$ws.Cells.Item(16, 1).Value = "Neutrophils"
$ws.Cells.Item(16, 2).Value = "Ntng1"
$ws.Cells.Item(16, 3).Value = "Lrrc4c"
$ws.Cells.Item(16, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.002167
$ws.Cells.Item(16, 8).Value = 0.006501
$ws.Cells.Item(16, 9).Value = 0.01200760241371094
$ws.Cells.Item(16, 10).Value = 0.01200760241371094
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.0009266666666666667
$ws.Cells.Item(16, 14).Value = 0.00278
$ws.Cells.Item(16, 15).Value = 0.001086631925715027
$ws.Cells.Item(16, 16).Value = 0.001086631925715027
$ws.Cells.Item(16, 17).Value = 0.000002008086666666667
$ws.Cells.Item(16, 18).Value = 0.00001807278
$ws.Cells.Item(16, 19).Value = 0.00001304784413403113
$ws.Cells.Item(16, 20).Value = 0.00001304784413403113

# Row 17
$ws.Cells.Item(17, 1).Value = "Neutrophils"
$ws.Cells.Item(17, 2).Value = "Ntng1"
$ws.Cells.Item(17, 3).Value = "Lrrc4c"
$ws.Cells.Item(17, 4).Value = "MuSCs"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.002167
$ws.Cells.Item(17, 8).Value = 0.006501
$ws.Cells.Item(17, 9).Value = 0.01200760241371094
$ws.Cells.Item(17, 10).Value = 0.01200760241371094
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.5128903333333333
$ws.Cells.Item(17, 14).Value = 1.538671
$ws.Cells.Item(17, 15).Value = 0.6014277092704556
$ws.Cells.Item(17, 16).Value = 0.6014277092704556
$ws.Cells.Item(17, 17).Value = 0.001111433352333333
$ws.Cells.Item(17, 18).Value = 0.010002900171
$ws.Cells.Item(17, 19).Value = 0.007221704813508565
$ws.Cells.Item(17, 20).Value = 0.007221704813508567
